$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.952.44"
$ws.Range("E2").Value = "  -3.24%  "

# Row 3
$ws.Range("D3").Value = "2.565.78"
$ws.Range("E3").Value = "  -1.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.39%  "

# Row 5
$ws.Range("Z1").Formula = "=""505.61"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  -3.45%  "

# Row 6
$ws.Range("Z1").Formula = "=""144.74"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  -6.41%  "

# Row 7
$ws.Range("E7").Value = "  +0.52%  "

# Row 8
$ws.Range("Z1").Formula = "=""0.555"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Range("E8").Value = "  -6.44%  "

# Row 9
$ws.Range("D9").Value = "2.559.30"
$ws.Range("E9").Value = "  -1.75%  "

# Row 10
$ws.Range("Z1").Formula = "=""6.17"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = "  -7.87%  "

# Row 11
$ws.Range("Z1").Formula = "=""0.102"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = "  -3.46%  "

# Row 12
$ws.Range("Z1").Formula = "=""0.331"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  -4.90%  "

# Row 13
$ws.Range("E13").Value = "  -1.16%  "

# Row 14
$ws.Range("D14").Value = "3.029.09"
$ws.Range("E14").Value = "  -0.72%  "

# Row 15
$ws.Range("D15").Value = "59.244.65"
$ws.Range("E15").Value = "  -2.80%  "

# Row 16
$ws.Range("Z1").Formula = "=""20.60"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Range("E16").Value = "  -4.90%  "

# Row 17
$ws.Range("Z1").Formula = "=""0.0000134"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = "  -4.88%  "

# Row 18
$ws.Range("D18").Value = "2.582.83"
$ws.Range("E18").Value = "  -0.65%  "

# Row 19
$ws.Range("Z1").Formula = "=""4.52"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  -4.87%  "

# Row 20
$ws.Range("Z1").Formula = "=""334.50"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$ws.Range("E20").Value = "  -5.42%  "

# Row 21
$ws.Range("Z1").Formula = "=""10.10"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  -4.76%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("Z1").Formula = "=""0.996"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  -0.31%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("Z1").Formula = "=""5.94"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Value = "  -4.81%  "

# Row 24
$ws.Range("Z1").Formula = "=""60.07"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Range("E24").Value = "  -1.43%  "

# Row 25
$ws.Range("Z1").Formula = "=""0.408"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  -4.47%  "

# Row 26
$ws.Range("Z1").Formula = "=""1.01"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$ws.Range("E26").Value = "  +0.80%  "

# Row 27
$ws.Range("Z1").Formula = "=""0.155"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  -6.54%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0776"
$ws.Range("E28").Value = "  -8.22%  "

# Row 29
$ws.Range("Z1").Formula = "=""6.87"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  -7.05%  "

# Row 30
$ws.Range("E30").Value = "  +0.13%  "

# Row 31
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("Z1").Formula = "=""5.85"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  -7.65%  "

# Row 32
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("Z1").Formula = "=""149.14"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Value = "  +0.48%  "

# Row 33
$ws.Range("Z1").Formula = "=""18.54"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = "  -4.31%  "

# Row 34
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("Z1").Formula = "=""1.54"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  -3.74%  "

# Row 35
$ws.Range("Z1").Formula = "=""3.85"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  -7.64%  "

# Row 36
$ws.Range("Z1").Formula = "=""0.883"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  -6.13%  "

# Row 37
$ws.Range("Z1").Formula = "=""1.11"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  -7.87%  "

# Row 38
$ws.Range("Z1").Formula = "=""36.25"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  -0.64%  "

# Row 39
$ws.Range("Z1").Formula = "=""0.823"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Range("E39").Value = "  -3.25%  "

# Row 40
$ws.Range("Z1").Formula = "=""1.38"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  -8.14%  "

# Row 41
$ws.Range("Z1").Formula = "=""3.52"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = "  -7.36%  "

# Row 42
$ws.Range("Z1").Formula = "=""281.36"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  -2.57%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("Z1").Formula = "=""1.00"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  +0.31%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("Z1").Formula = "=""0.608"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Value = "  -2.53%  "

# Row 45
$ws.Range("Z1").Formula = "=""0.0980"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  -3.30%  "

# Row 46
$ws.Range("Z1").Formula = "=""0.0532"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = "  -5.04%  "

# Row 47
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("Z1").Formula = "=""10.34"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  +0.10%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("Z1").Formula = "=""18.61"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  -5.15%  "

# Row 49
$ws.Range("Z1").Formula = "=""0.0228"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  -4.41%  "

# Row 50
$ws.Range("Z1").Formula = "=""4.51"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  -7.75%  "

# Row 51
$ws.Range("D51").Value = "1.910.87"
$ws.Range("E51").Value = "  -2.40%  "

$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = $false
